$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round Q7 and R7 to whole numbers
$ws.Range("Q7").Value = 715208
$ws.Range("R7").Value = 7303647

# Clear Z7 and AB7 (Starttid / Sluttid) cell contents
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
